$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ50203817",
    "summ50388717",
    "summ50556884",
    "summ50742797",
    "summ50920815",
    "summ51104823",
    "summ51272215",
    "summ51455325",
    "summ51638105",
    "summ51804901",
    "summ52003546",
    "summ52172323",
    "summ52355971",
    "summ52540313",
    "summ52723398",
    "summ52907344",
    "summ53090151",
    "summ53272525",
    "summ53455867",
    "summ53737274",
    "summ53914241",
    "summ54088958",
    "summ54278325",
    "summ54459967",
    "summ54639343",
    "summ54825627",
    "summ55035054",
    "summ55240873",
    "summ55420342",
    "summ55591351",
    "summ55775379",
    "summ55958636",
    "summ56138034",
    "summ56326349",
    "summ56507368",
    "summ56689665",
    "summ56887000",
    "summ57073690",
    "summ57259700",
    "summ57441859",
    "summ57622637",
    "summ57827595",
    "summ58035371",
    "summ58230303",
    "summ58405475",
    "summ58588919",
    "summ58771991",
    "summ58955657",
    "summ59138236",
    "summ59321928"
)

$newValues = @(
    @("2348.305469151413", "0.009693115658345901", "417.5531579586019", "1.306059574977568e-06", "-0.09546449873587767", "0.2814249302134852", "-3059.822611456061", "0.3558231088910718", "1171.083205950103", "0.2857668138759802", "307.2535883705019", "0.007836658172734949"),
    @("1692.828946716686", "0.1000483883038616", "476.5860258254478", "8.478142447121008e-06", "-0.05910690490414316", "0.5906153356977799", "-170.6305286844963", "0.9620949980285691", "707.6427624406231", "0.621208015492986", "333.2879189944538", "0.02508145903068"),
    @("1001.691930359139", "0.2502963498474452", "516.3040521495413", "5.82913335215281e-07", "-0.07149521425161323", "0.4676591712052658", "2412.919147276549", "0.4737359674649626", "1622.288288967044", "0.2215767238625535", "287.9778687136451", "0.01984068516654152"),
    @("948.0885096254488", "0.3601124230800636", "494.918581504507", "2.295738261328439e-06", "-0.06964049597163607", "0.5433904648912091", "3910.391157965183", "0.319714589487498", "1048.098178618945", "0.4354848822302583", "296.0437339101115", "0.02652021497753916"),
    @("1780.251067973491", "0.08761558741551498", "422.1243812555658", "5.623871082550575e-05", "-0.06672918261701495", "0.5980574654909203", "1437.517069916298", "0.6786158828930471", "722.7655219575959", "0.5947321090742086", "295.8113012585528", "0.02891795152647609"),
    @("1826.269634203486", "0.03836784486698702", "435.1976661260045", "1.171571583401333e-06", "-0.05917846215662326", "0.5374851463797498", "-1242.956283517469", "0.6767050376795194", "860.9461175673719", "0.4245209869305034", "347.3027028625634", "0.006584120607595599"),
    @("1789.284133798462", "0.06146065687889164", "476.6077984397691", "1.682397056099143e-06", "-0.06370353889707245", "0.5675358756723639", "2142.705190395976", "0.5593606037105481", "147.9318977061689", "0.9241675473252569", "224.4766227883296", "0.08360809784437527"),
    @("1265.491031930064", "0.1815173414255524", "520.1733642542786", "3.568034778100543e-06", "-0.05954707912767687", "0.6128638040151146", "804.3259268288516", "0.8141559290297242", "1486.499383636455", "0.2362098200462069", "273.8083988972878", "0.02508976950240157"),
    @("983.1160651758478", "0.3581347369315246", "467.3012598772439", "2.655573490214377e-05", "-0.0446206146831425", "0.6729921040728902", "3816.725939036511", "0.3339155193061656", "712.0812773373946", "0.5862751957801235", "334.9965336576677", "0.01938349806949894"),
    @("1776.684681777462", "0.09833873354569987", "441.9933549649526", "1.596720930091483e-05", "-0.1060574736642687", "0.3381261319254094", "342.6470145428357", "0.927050693486531", "1491.368281798256", "0.3032774700385963", "307.206216957462", "0.02759656916110139"),
    @("1484.541761303143", "0.08325444286386506", "452.7263685364773", "2.594957761672882e-07", "-0.04223307558492723", "0.6063704136962433", "342.1560703015639", "0.9091359831332591", "913.7349775780224", "0.371380165606484", "324.2792914234214", "0.007616740106452843"),
    @("1341.681957511726", "0.1309051064985562", "401.2872316910949", "4.461560279086112e-05", "-0.08119882970160919", "0.3984961717692682", "5519.12995510242", "0.1257046981779499", "-60.69390583131803", "0.9599361087511803", "349.5092016195103", "0.0121174922762121"),
    @("1891.445586040454", "0.05880706298546392", "452.1482618854143", "2.404835453383261e-05", "-0.103635000215101", "0.332615712164663", "368.356579284201", "0.9198739681317076", "1331.384614514015", "0.3313747840950018", "279.0242813619595", "0.0499045576557559"),
    @("1224.852419995563", "0.2015585814579073", "567.5574810551925", "6.060178500186762e-07", "-0.03159606814136373", "0.758108464157764", "-673.323625921772", "0.8505278896242523", "1403.638201437702", "0.2498076426081781", "261.3178366717715", "0.0259074067614948"),
    @("1645.275351395942", "0.1464155353179784", "461.6733201314473", "7.829304579289486e-06", "-0.101411945403727", "0.6290943785312996", "274.724696061352", "0.9428880402728621", "1360.803536619126", "0.5427416135715857", "311.024457286099", "0.02582467230683121"),
    @("1022.560721078869", "0.2438253827849603", "587.4250747135666", "1.216518477014734e-07", "-0.04826667999427858", "0.6287996721959901", "129.1290379749171", "0.9667550722867421", "1982.244016560518", "0.1290198272456058", "268.5603506370817", "0.02219781569975135"),
    @("1871.999925312244", "0.0943685657140881", "447.5855323045753", "1.946309081387114e-05", "-0.07581048898883225", "0.6543780714792278", "495.6684722205264", "0.8976877486886581", "874.2277204290394", "0.6074675463591743", "280.7992623560428", "0.04886447448491163"),
    @("2051.221946389231", "0.02927126708148018", "380.6698582107207", "9.407744974713255e-05", "-0.08601178546532598", "0.3807897889405674", "2596.781495010498", "0.4465917477339866", "38.46128960548231", "0.9759632148868014", "261.0796801619506", "0.04205427373499341"),
    @("1376.412968261905", "0.1642427630321372", "424.5512077300718", "3.313267174212258e-05", "-0.09626641042353637", "0.3740561214911099", "3336.665082774656", "0.4542769558529229", "749.56981430617", "0.5704870870635236", "371.0527657471833", "0.0100427573565312"),
    @("1249.134719628049", "0.1614764316470568", "483.5840224394481", "2.615999011284445e-08", "-0.06198826608790906", "0.4276148215515215", "-245.2263317844463", "0.9333519255379037", "1445.182305810446", "0.1442152542126593", "336.2285158939162", "0.003455627886541806"),
    @("1030.327441887664", "0.2756824899404758", "546.4954594319649", "9.004131830345794e-08", "-0.02696823734585271", "0.7624229258454802", "3929.770113310676", "0.2551434825915169", "344.0426343780359", "0.7653326378243495", "206.1570780449429", "0.1048841202118766"),
    @("1590.137599709236", "0.1093000449711331", "430.6885444450398", "5.173715603808959e-05", "-0.1054736931198519", "0.3875337936048998", "2144.460469906369", "0.5628271886826727", "1291.746444009021", "0.3539288005265989", "295.3471391033765", "0.04346709711353926"),
    @("1905.917270300928", "0.07523871896800841", "434.0768099807752", "4.062197892080797e-05", "-0.06806140362921903", "0.597799284466205", "-513.9667653294791", "0.8964367295541917", "840.1822576626655", "0.5590678377656486", "315.1551028866849", "0.01682572288087359"),
    @("1524.536267648481", "0.06448880138887404", "484.468500450662", "6.97520825450899e-08", "-0.0758598823226902", "0.3713493679564845", "-597.5377849993884", "0.8348972731778608", "1837.754372181179", "0.1208900921374937", "307.809033172137", "0.007856078004056759"),
    @("1923.593691725416", "0.0690416313472394", "439.9286371432586", "1.724679015613266e-05", "-0.08767895152987965", "0.4256246930438441", "-549.4559159014934", "0.884525123820759", "841.4405556844301", "0.5473678269256085", "331.110559069338", "0.0164382945072443"),
    @("1777.292376917", "0.08917389529383772", "466.2752117166405", "1.527200501897081e-05", "-0.05126199909728085", "0.62687451138725", "434.9464214976215", "0.9065992586947282", "650.2647339506457", "0.6151296230676314", "271.0861439812376", "0.05139149653327254"),
    @("903.1788339393338", "0.3782137636803773", "474.0696818678252", "3.616141895514974e-06", "-0.06829649004223526", "0.5799018670695896", "4080.517187177118", "0.2946322364420966", "988.8230183229948", "0.4964893997054018", "326.4626362844622", "0.01361125060118861"),
    @("1845.450442709798", "0.08086890701508251", "450.1304747395367", "4.254391030026509e-05", "-0.1309507005113836", "0.2374381007548968", "39.08763088180604", "0.9913357267228241", "1655.324172922544", "0.227289869165361", "303.8184751850634", "0.03463507093527032"),
    @("1588.752587407877", "0.01969639825929849", "481.2877004553617", "6.041682710724428e-09", "-0.02529383164567267", "0.7098789407680945", "153.0427650145921", "0.9535214039089757", "667.1664138876922", "0.4602077872572325", "261.4879900761615", "0.0050710531901062"),
    @("1794.342921599826", "0.1110203600321677", "433.1387152304907", "4.957219407687026e-05", "-0.07429089553301368", "0.5790889565870756", "-48.65877846661169", "0.9901147288956846", "951.578632934049", "0.5447927939523558", "351.8717432556469", "0.01878103396024337"),
    @("1275.34160475169", "0.2843415975728908", "505.1648472039599", "9.666457293794635e-07", "-0.08257429247623277", "0.4852378494330279", "1624.614332419931", "0.6766116289712777", "1197.762396484432", "0.3723213714610176", "270.8970639442985", "0.05858019955259099"),
    @("1894.6050676242", "0.02005856141788255", "442.9305041096151", "6.441586051863229e-07", "-0.04946865228416755", "0.5774188472205193", "-2266.261527971218", "0.4338599119219227", "975.0515632851852", "0.3634979286840095", "348.7269733482885", "0.001264681378945932"),
    @("1313.937234200929", "0.1664490714752173", "466.1327680722465", "2.155210246033626e-05", "-0.1093847655074242", "0.2646588528955889", "1598.537570259437", "0.6023106809796304", "2144.498693296234", "0.09791700957915457", "314.3203292910781", "0.0195340061252566"),
    @("1548.72957143149", "0.1413808435490934", "434.6921658032143", "2.857743439685555e-05", "-0.07556807924353198", "0.5802622063431249", "2325.10115196671", "0.5947413403854878", "700.1608594738891", "0.6600912933608101", "326.4566948151281", "0.02605151208652062"),
    @("1499.913644440716", "0.123376018268543", "477.5432987384196", "2.746364578811763e-06", "-0.05273790993228575", "0.5940935381309249", "2998.428489025941", "0.428770591272894", "204.9676174354636", "0.8747666858197719", "250.435917549727", "0.0567792882121213"),
    @("1567.257647184987", "0.1486301784721042", "473.6286538672628", "3.993839501054015e-06", "-0.0550537767761663", "0.578046926731542", "2002.159470667419", "0.5775828774019656", "299.1430757219437", "0.8203849962772189", "281.2598172268313", "0.03787618762328628"),
    @("1460.83798092393", "0.1339494736271515", "357.0951164983296", "0.0007125804510379435", "-0.1479936235671799", "0.1853273557809822", "3265.658968265103", "0.440920190480838", "1433.248822309455", "0.3033924304734475", "458.2212411248246", "0.006148649448305218"),
    @("1431.363871384039", "0.2005667995274856", "483.7369446260106", "1.966922014812871e-06", "-0.139667623375426", "0.4317280872627298", "965.5109976706844", "0.8038196350142013", "1618.302616463999", "0.2825798530580017", "291.4535753856994", "0.02332626538717071"),
    @("2405.350597542383", "0.002922595288576726", "380.2973149891534", "1.104491681564876e-05", "-0.09692300118714942", "0.2499111197176229", "-1788.33562952069", "0.5197836795942687", "1025.592041114641", "0.3166075458362533", "313.1948678437183", "0.003767229018687244"),
    @("863.3453221180764", "0.3610505363093758", "584.4387160243911", "9.347182618930343e-08", "-0.02647279377995876", "0.8010059946256236", "1521.099068766991", "0.6406865208881556", "1704.961573406558", "0.2356882984556882", "235.0105009107518", "0.06195622746498355"),
    @("1448.297257396999", "0.0224324976954014", "521.3715132547572", "1.387061735473329e-09", "-0.03337810855017025", "0.6298058675481999", "456.0888280155996", "0.8417674962684548", "974.6916606418738", "0.3290208124170544", "238.6903852469143", "0.006486398253090606"),
    @("1591.010314595396", "0.1314221914863181", "458.5547414365849", "1.541850775033755e-05", "-0.09095291250091664", "0.5788710405323383", "1264.260935228132", "0.7428346050761372", "922.1006070009407", "0.5301812668417619", "311.8947966608534", "0.03606787578623948"),
    @("2006.565348138477", "0.08302849029951925", "468.5985203593756", "9.632722772378845e-06", "-0.05831469810252753", "0.5833652063052304", "115.3385231248731", "0.9762316484827869", "560.0271174318", "0.6784575515212532", "253.2164001385497", "0.08468969794829095"),
    @("2095.176344314988", "0.04420246803043022", "389.9263242239981", "0.0001217284949673653", "-0.1738404568954923", "0.146247200891378", "65.06022535817283", "0.9860908648735333", "2039.778022419972", "0.1769416219109232", "322.7579397974225", "0.02213072096341891"),
    @("597.6389327605941", "0.5661992654547218", "482.0358594375765", "3.88896805956668e-06", "-0.03802362176940341", "0.7035157327460426", "3602.54466641549", "0.3608505373893498", "965.5086561326366", "0.4292900041746855", "380.4360113024773", "0.004744910990269089"),
    @("2176.487266560496", "0.04089848072441597", "401.2888678988716", "0.0001577958015205151", "-0.0835316061487815", "0.4403236559012346", "-771.4703087051766", "0.8353234871156157", "604.9290783488423", "0.6629912010538835", "342.0413590626603", "0.01857812592031461"),
    @("1777.183800241304", "0.1222684614455936", "470.0208525147991", "1.302644939297991e-05", "-0.06431609293225715", "0.5715744289633262", "-117.8087693038808", "0.9755755621032178", "970.6495062882418", "0.4744751603677412", "275.9713667658528", "0.04894234456615169"),
    @("1325.36354301398", "0.1258879674440475", "435.7977393645446", "2.134425560838362e-07", "-0.0007248730112526247", "0.995224448718018", "337.8258019721106", "0.9144389026809602", "998.237361061264", "0.3939758865250611", "354.3511125458417", "0.001206335236492845"),
    @("1029.932821889206", "0.2294258702971975", "578.5750035663972", "8.312397120189833e-08", "-0.03037403123376356", "0.748172812621425", "1663.192058753379", "0.610163728335688", "1058.565797103117", "0.4282066710376753", "242.6278484825239", "0.04969358904613942"),
    @("1456.026371805453", "0.1475587133342055", "440.7928656013612", "2.15148245933303e-05", "-0.1143860629341676", "0.2741791777331917", "3209.573121885624", "0.3903177319296796", "1170.106657663914", "0.3649297041131816", "306.7275872452353", "0.03141076720374934")
)

for ($i = 0; $i -lt 50; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
    $vals = $newValues[$i]
    $row = 2
    for ($j = 0; $j -lt 12; $j += 2) {
        $ws.Cells.Item($row, 2).Value = [double]$vals[$j]
        $ws.Cells.Item($row, 3).Value = [double]$vals[$j + 1]
        $row += 1
    }
}

Write-Output "Applied updates to all sheets"